$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Sheet "samples_retained": update the Emotional_EMA row (row 13).
# Replace the old formulas in C13/D13 with plain values, and update the
# notes cell (H13) with the new note text.
# ---------------------------------------------------------------------------
$wsRetained = $wb.Worksheets.Item("samples_retained")
$wsRetained.Range("H13").Value = "perceived valence recoded to majority and average votes; if either matched intended valence, sample kept; discarded 58"
$wsRetained.Range("C13").Value = 147
$wsRetained.Range("D13").Value = 333

# ---------------------------------------------------------------------------
# Sheet "positive": add the new "interest"/"intrest" rows (reuse the "int"
# abbreviation) for the reworked Emotional_EMA data organisation.
# ---------------------------------------------------------------------------
$wsPositive = $wb.Worksheets.Item("positive")
$wsPositive.Range("A6").Value = "interest"
$wsPositive.Range("B6").Value = "en"
$wsPositive.Range("C6").Value = "int"

$wsPositive.Range("A7").Value = "intrest"
$wsPositive.Range("B7").Value = "en"
$wsPositive.Range("C7").Value = "int"

# ---------------------------------------------------------------------------
# Sheet "negative": add the new emotion/abbreviation rows for Emotional_EMA.
# ---------------------------------------------------------------------------
$wsNegative = $wb.Worksheets.Item("negative")
$wsNegative.Range("A19").Value = "Langeweile"
$wsNegative.Range("B19").Value = "de"
$wsNegative.Range("C19").Value = "bor"

$wsNegative.Range("A20").Value = "boredom"
$wsNegative.Range("B20").Value = "en"
$wsNegative.Range("C20").Value = "bor"

$wsNegative.Range("A21").Value = "surprise"
$wsNegative.Range("B21").Value = "en"
$wsNegative.Range("C21").Value = "sur"

$wsNegative.Range("A22").Value = "surprised"
$wsNegative.Range("B22").Value = "en"
$wsNegative.Range("C22").Value = "sur"

$wsNegative.Range("A23").Value = "surprisse"
$wsNegative.Range("B23").Value = "fr"
$wsNegative.Range("C23").Value = "sur"

$wsNegative.Range("A24").Value = "unsure"
$wsNegative.Range("B24").Value = "en"
$wsNegative.Range("C24").Value = "unc"

$wsNegative.Range("A25").Value = "uncertain"
$wsNegative.Range("B25").Value = "en"
$wsNegative.Range("C25").Value = "unc"

$wsNegative.Range("A26").Value = "frustration"
$wsNegative.Range("B26").Value = "en"
$wsNegative.Range("C26").Value = "fru"

# ---------------------------------------------------------------------------
# Sheet "positive": add the remaining "excitement"/"curiosity" rows.
# ---------------------------------------------------------------------------
$wsPositive.Range("A8").Value = "excitement"
$wsPositive.Range("B8").Value = "en"
$wsPositive.Range("C8").Value = "exc"

$wsPositive.Range("A9").Value = "curiosity"
$wsPositive.Range("B9").Value = "en"
$wsPositive.Range("C9").Value = "cur"

# ---------------------------------------------------------------------------
# Sheet "discard": several of the old rows have been reclassified into the
# "positive" / "negative" sheets above, so remove them here and leave only
# the genuinely leftover/ambiguous terms (amused, bothered, concentrating,
# sleepiness, sleepy, thinking).
# ---------------------------------------------------------------------------
$wsDiscard = $wb.Worksheets.Item("discard")
$wsDiscard.Range("A1:B14").ClearContents()

$wsDiscard.Range("A1").Value = "amused"
$wsDiscard.Range("B1").Value = "en"

$wsDiscard.Range("A2").Value = "bothered"
$wsDiscard.Range("B2").Value = "en"

$wsDiscard.Range("A3").Value = "concentrating"
$wsDiscard.Range("B3").Value = "en"

$wsDiscard.Range("A4").Value = "sleepiness"
$wsDiscard.Range("B4").Value = "en"

$wsDiscard.Range("A5").Value = "sleepy"
$wsDiscard.Range("B5").Value = "en"

$wsDiscard.Range("A6").Value = "thinking"
$wsDiscard.Range("B6").Value = "en"
